$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 150, pushing existing rows 150-207 down to 151-208.
# Excel's native row-insert inherits formatting (incl. the date NumberFormat on
# column D) from the row above, matching the rest of the data block.
$ws.Rows("150:150").Insert()

# Populate the newly inserted (now empty) row 150 with the new record.
$ws.Range("A150").Value2 = 10
$ws.Range("B150").Value2 = "Vega Modelo de Temuco"
$ws.Range("C150").Value2 = "La Araucanía"
$ws.Range("D150").Value2 = 44704
$ws.Range("E150").Value2 = 9
$ws.Range("F150").Value2 = 100112005
$ws.Range("G150").Value2 = "Puerro"
$ws.Range("H150").Value2 = "Azul de Maquehue"
$ws.Range("I150").Value2 = "Primera"
$ws.Range("J150").Value2 = 30
$ws.Range("K150").Value2 = 12000
$ws.Range("L150").Value2 = 12000
$ws.Range("M150").Value2 = 12000
$ws.Range("N150").Value2 = "$/docena de paquetes"
$ws.Range("O150").Value2 = "Provincia de Cautín"
$ws.Range("P150").Value2 = 1000
$ws.Range("Q150").Value2 = 12
$ws.Range("R150").Value2 = "Hortaliza"
